$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1241.25
$ws.Range("I53").Value = 1676.5
$ws.Range("K53").Value = 1676.5
$ws.Range("M53").Value = -1039.5
$ws.Range("H74").Value = 3777.077
$ws.Range("I74").Value = 3425.25
$ws.Range("K74").Value = 3425.25
$ws.Range("M74").Value = -2489.25
$ws.Range("H76").Value = 3999.25
$ws.Range("I76").Value = 3999.1
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3999.1
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3684.1
$ws.Range("N76").Value = -4630
$ws.Range("H77").Value = 3777.077
$ws.Range("I77").Value = 3425.25
$ws.Range("K77").Value = 17126.25
$ws.Range("M77").Value = -12446.25
$ws.Range("H79").Value = 3999.25
$ws.Range("I79").Value = 3999.1
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3999.1
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2907.1
$ws.Range("N79").Value = -6184
$ws.Range("H92").Value = 677.4
$ws.Range("I92").Value = 648.625
$ws.Range("K92").Value = 648.625
$ws.Range("M92").Value = 599.375
$ws.Range("H103").Value = 2761.7144
$ws.Range("I103").Value = 3597.4
$ws.Range("K103").Value = 10792.2
$ws.Range("M103").Value = -10206.2
$ws.Range("H127").Value = 1185
$ws.Range("I127").Value = 1185
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 3555
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1405
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 2295.5957
$ws.Range("I129").Value = 491
$ws.Range("J129").Value = 2418.6365
$ws.Range("K129").Value = 1473
$ws.Range("L129").Value = 7255.9095
$ws.Range("M129").Value = 3527
$ws.Range("N129").Value = -17255.9095
$ws.Range("H132").Value = 35778.434
$ws.Range("I132").Value = 36908.723
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 110726.169
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -108196.169
$ws.Range("N132").Value = -14060
$ws.Range("H141").Value = 5041.7144
$ws.Range("I141").Value = 4499.6665
$ws.Range("J141").Value = 5448.25
$ws.Range("K141").Value = 13498.9995
$ws.Range("L141").Value = 16344.75
$ws.Range("M141").Value = -8318.999500000002
$ws.Range("N141").Value = -26704.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2973.5757
$ws.Range("I32").Value = 3037.6453
$ws.Range("J32").Value = 1980.5
$ws.Range("K32").Value = 3037.6453
$ws.Range("L32").Value = 1980.5
$ws.Range("M32").Value = -2750.6453
$ws.Range("N32").Value = -2554.5
$ws.Range("H45").Value = 1500.2174
$ws.Range("I45").Value = 1440.1666
$ws.Range("K45").Value = 1440.1666
$ws.Range("M45").Value = -1063.1666
$ws.Range("H102").Value = 1175.7273
$ws.Range("I102").Value = 1093.3
$ws.Range("K102").Value = 1093.3
$ws.Range("M102").Value = 528.7
$ws.Range("H122").Value = 1434.6316
$ws.Range("I122").Value = 985.5714
$ws.Range("K122").Value = 2956.7142
$ws.Range("M122").Value = -506.7142000000003
$ws.Range("H124").Value = 53451.668
$ws.Range("J124").Value = 53451.668
$ws.Range("L124").Value = 53451.668
$ws.Range("N124").Value = -63271.668
$ws.Range("H125").Value = 84444
$ws.Range("J125").Value = 84444
$ws.Range("L125").Value = 84444
$ws.Range("N125").Value = -94284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2359.5557
$ws.Range("I105").Value = 2185.75
$ws.Range("K105").Value = 2185.75
$ws.Range("M105").Value = -438.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2552.2273
$ws.Range("I31").Value = 1872.3636
$ws.Range("K31").Value = 1872.3636
$ws.Range("M31").Value = -1577.3636
$ws.Range("H34").Value = 2552.2273
$ws.Range("I34").Value = 1872.3636
$ws.Range("K34").Value = 1872.3636
$ws.Range("M34").Value = -1670.3636
$ws.Range("H58").Value = 2434.5
$ws.Range("I58").Value = 2434.5
$ws.Range("K58").Value = 2434.5
$ws.Range("M58").Value = -2231.5
$ws.Range("H97").Value = 44844.5
$ws.Range("J97").Value = 44844.5
$ws.Range("L97").Value = 44844.5
$ws.Range("N97").Value = -46826.5
$ws.Range("H99").Value = 2176.4285
$ws.Range("I99").Value = 2303.5
$ws.Range("K99").Value = 2303.5
$ws.Range("M99").Value = -805.5
$ws.Range("H126").Value = 2176.4285
$ws.Range("I126").Value = 2303.5
$ws.Range("K126").Value = 6910.5
$ws.Range("M126").Value = -4440.5
$ws.Range("H132").Value = 3100
$ws.Range("I132").Value = 2200
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6600
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4070
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 2434.5
$ws.Range("I136").Value = 2434.5
$ws.Range("K136").Value = 7303.5
$ws.Range("M136").Value = -4753.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 505
$ws.Range("J12").Value = 618.6667
$ws.Range("L12").Value = 1856.0001
$ws.Range("N12").Value = -2202.0001
$ws.Range("H97").Value = 134109
$ws.Range("I97").Value = 159819.8
$ws.Range("J97").Value = 5555
$ws.Range("K97").Value = 479459.4
$ws.Range("L97").Value = 16665
$ws.Range("M97").Value = -478963.4
$ws.Range("N97").Value = -17657
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H131").Value = 869.8889
$ws.Range("I131").Value = 599.5
$ws.Range("J131").Value = 3033
$ws.Range("K131").Value = 1798.5
$ws.Range("L131").Value = 9099
$ws.Range("M131").Value = 3241.5
$ws.Range("N131").Value = -19179

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 18303.4
$ws.Range("J35").Value = 18412.5
$ws.Range("L35").Value = 18412.5
$ws.Range("N35").Value = -19008.5
$ws.Range("H80").Value = 3011.6
$ws.Range("I80").Value = 3249.5
$ws.Range("J80").Value = 2060
$ws.Range("K80").Value = 3249.5
$ws.Range("L80").Value = 2060
$ws.Range("M80").Value = -2251.5
$ws.Range("N80").Value = -4056
$ws.Range("H83").Value = 3011.6
$ws.Range("I83").Value = 3249.5
$ws.Range("J83").Value = 2060
$ws.Range("K83").Value = 16247.5
$ws.Range("L83").Value = 10300
$ws.Range("M83").Value = -11255.5
$ws.Range("N83").Value = -20284
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H107").Value = 753.53845
$ws.Range("I107").Value = 806.13043
$ws.Range("K107").Value = 806.13043
$ws.Range("M107").Value = 1113.86957

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 765.2
$ws.Range("J16").Value = 989.5
$ws.Range("L16").Value = 989.5
$ws.Range("N16").Value = -1329.5
$ws.Range("H68").Value = 4756.0586
$ws.Range("I68").Value = 2583.375
$ws.Range("K68").Value = 2583.375
$ws.Range("M68").Value = -1834.375
$ws.Range("H71").Value = 4756.0586
$ws.Range("I71").Value = 2583.375
$ws.Range("K71").Value = 12916.875
$ws.Range("M71").Value = -9172.875
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 7303
$ws.Range("I93").Value = 7303
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 7303
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -6055
$ws.Range("N93").ClearContents()
$ws.Range("H136").Value = 3444.182
$ws.Range("I136").Value = 1998.625
$ws.Range("K136").Value = 5995.875
$ws.Range("M136").Value = -3445.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7811.85
$ws.Range("I81").Value = 9269.532999999999
$ws.Range("J81").Value = 3438.8
$ws.Range("K81").Value = 18539.066
$ws.Range("L81").Value = 6877.6
$ws.Range("M81").Value = -17478.066
$ws.Range("N81").Value = -8999.6
$ws.Range("H84").Value = 7811.85
$ws.Range("I84").Value = 9269.532999999999
$ws.Range("J84").Value = 3438.8
$ws.Range("K84").Value = 92695.32999999999
$ws.Range("L84").Value = 34388
$ws.Range("M84").Value = -87391.32999999999
$ws.Range("N84").Value = -44996
$ws.Range("H100").Value = 1540.5625
$ws.Range("I100").Value = 1122
$ws.Range("J100").Value = 1959.125
$ws.Range("K100").Value = 2244
$ws.Range("L100").Value = 3918.25
$ws.Range("M100").Value = -1703
$ws.Range("N100").Value = -5000.25
$ws.Range("H112").Value = 42500
$ws.Range("J112").Value = 42500
$ws.Range("L112").Value = 42500
$ws.Range("N112").Value = -45454
$ws.Range("H122").Value = 1649
$ws.Range("I122").Value = 1658.5333
$ws.Range("J122").Value = 1601.3334
$ws.Range("K122").Value = 4975.5999
$ws.Range("L122").Value = 4804.0002
$ws.Range("M122").Value = -2525.5999
$ws.Range("N122").Value = -9704.0002
$ws.Range("H132").Value = 1815.8889
$ws.Range("I132").Value = 1742.875
$ws.Range("K132").Value = 5228.625
$ws.Range("M132").Value = -2698.625

